$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain decimal numbers need to be forced to
# Text format first, otherwise Excel will auto-convert the string into a number
# (e.g. "0.9999" -> 0.9999) and the literal formatting (trailing zeros, etc.)
# would be lost.
$numericLooking = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D22","D23","D25","D26","D27","D29","D30","D31","D32","D33","D35","D36","D40","D41","D42","D43","D44","D45","D50")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptocurrency data
$ws.Range('D2').Value = '29.385.49'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.876.08'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '0.7118'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '242.06'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.3119'
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('D9').Value = '0.07795'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '25.17'
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').Value = '0.08464'
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').Value = '1.868.20'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = '5.235'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '0.7125'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '91.28'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '29.385.74'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '6.059'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').Value = '0.000008231'
$ws.Range('E18').Value = '  +5.31%  '
$ws.Range('D19').Value = '241.09'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').Value = '13.25'
$ws.Range('E20').Value = '  +0.69%  '
$ws.Range('D21').Value = '2.117.57'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = '0.9996'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = '7.784'
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').Value = '0.1597'
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '163.38'
$ws.Range('D27').Value = '9.072'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').Value = '1.512'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('D30').Value = '4.433'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '4.321'
$ws.Range('E31').Value = '  +1.56%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.287'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('D33').Value = '0.05291'
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = '1.180'
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '0.7451'
$ws.Range('E36').Value = '  -11.34%  '
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').Value = '1.228.76'
$ws.Range('E39').Value = '  +5.70%  '
$ws.Range('D40').Value = '2.723'
$ws.Range('E40').Value = '  +1.15%  '
$ws.Range('D41').Value = '6.481'
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8932'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '110.21'
$ws.Range('E43').Value = '  +8.04%  '
$ws.Range('D44').Value = '72.77'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '0.9997'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '2.015.26'
$ws.Range('E46').Value = '  -0.56%  '
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('E49').Value = '  +4.41%  '
$ws.Range('D50').Value = '9.388'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('E51').Value = '  +1.33%  '
